# Apply localization fix: replace curly double quotes (" ") used around
# quoted terms in the English (en_US) column with straight single quotes,
# per commit "update on 20210731 画中人".
# Curly apostrophes within words (e.g. I've, You're, couldn't) are left intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$LDQ = [char]0x201C   # “
$RDQ = [char]0x201D   # ”
$SQ  = "'"

function Fix-Quotes([string]$text) {
    return $text.Replace($LDQ, $SQ).Replace($RDQ, $SQ)
}

$cells = @("C70", "C88", "C94", "C95", "C98", "C124", "C129")

foreach ($cellRef in $cells) {
    $range = $ws.Range($cellRef)
    $current = [string]$range.Text
    $updated = Fix-Quotes($current)
    $range.Value2 = $updated
}

$wb.Save()
